# Update gh-pages to output generated at 456a3b4
# Sheet 1 = 展览 (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 465
$ws1.Range("F4").Value = 7781
$ws1.Range("F8").Value = 27
$ws1.Range("F17").Value = 5642
$ws1.Range("F18").Value = 160
$ws1.Range("F20").Value = 1066
$ws1.Range("F22").Value = 333

# Sheet 4 = 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 465
$ws4.Range("F4").Value = 7781
$ws4.Range("F8").Value = 27
$ws4.Range("F18").Value = 5642
$ws4.Range("F20").Value = 160
$ws4.Range("F22").Value = 1066
$ws4.Range("F24").Value = 333
